# Web server integration added
#
# Adds a new "GingerBreadMan" tag entry (with two tag codes) to the
# MakerRangerTags sheet. This is implemented as inserting two new rows
# just above the existing blank separator row that precedes the
# "Reprint/AbortGame/Player/Rounds/Show" block, which pushes that whole
# block down by two rows (old rows 32-37 -> new rows 34-39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 32 (existing rows 32-37 shift down to 34-39)
$ws.Rows("32:33").Insert()

# Fill row 33 first so its tag code (D7BE2FF4) is registered in the shared
# string table before row 32's tag code (A7C683F4), matching tag entry order.
$ws.Cells.Item(33, 1).Value = "GingerBreadMan"
$ws.Cells.Item(33, 2).Value = "10"
$ws.Cells.Item(33, 3).Value = "D7BE2FF4"
$ws.Cells.Item(33, 4).Formula = '=CONCATENATE(C33,CHAR(9),B33,"|",A33)'
$ws.Cells.Item(33, 5).Formula = '=LEFT(A33,16)'

$ws.Cells.Item(32, 1).Value = "GingerBreadMan"
$ws.Cells.Item(32, 2).Value = "10"
$ws.Cells.Item(32, 3).Value = "A7C683F4"
$ws.Cells.Item(32, 4).Formula = '=CONCATENATE(C32,CHAR(9),B32,"|",A32)'
$ws.Cells.Item(32, 5).Formula = '=LEFT(A32,16)'

# Match the "text" number format used by every other B/C cell in the table
$ws.Range("B32:C33").NumberFormat = "@"

# Update the view: scroll back to the top and select the whole D column
# range used by the data (D1:D39), as in the saved workbook state.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$ws.Range("D1:D39").Select()
